# Fix Training Data Issue (#48)
# The "Date" column (BF) for this sheet was populated with the file's own
# name/label ("2-18-2007-08") instead of an actual date string. Correct it
# to the real ISO-style date "2008-02-18" for every data row (rows 2-31).
#
# NOTE: Typing "2008-02-18" directly into a cell makes Excel auto-recognize
# it as a date and reformat/convert the cell. To keep the value as literal
# text (matching the original inline-string "Date" column), we enter the
# value with a leading apostrophe (forces text entry / quotePrefix) and
# then reapply the default "Normal" style so the cell is left without any
# extra number-formatting baggage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$dateCol = 58   # column BF

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $dateCol)
    $cell.Value = "'2008-02-18"
    $cell.Style = "Normal"
}
